$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2400
$ws.Range("I40").Value = 2100
$ws.Range("K40").Value = 2100
$ws.Range("M40").Value = -1925
$ws.Range("H100").Value = 1211.0625
$ws.Range("I100").Value = 958.0833
$ws.Range("J100").Value = 1970
$ws.Range("K100").Value = 958.0833
$ws.Range("L100").Value = 1970
$ws.Range("M100").Value = -417.0833
$ws.Range("N100").Value = -3052
$ws.Range("H141").Value = 793.9286
$ws.Range("I141").Value = 792.9167
$ws.Range("J141").Value = 800
$ws.Range("K141").Value = 2378.7501
$ws.Range("L141").Value = 2400
$ws.Range("M141").Value = 2801.2499
$ws.Range("N141").Value = -12760

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21282012
$ws.Range("I32").Value = 23813132
$ws.Range("J32").Value = 20602.6
$ws.Range("K32").Value = 23813132
$ws.Range("L32").Value = 20602.6
$ws.Range("M32").Value = -23812845
$ws.Range("N32").Value = -21176.6
$ws.Range("H45").Value = 1531.6
$ws.Range("I45").Value = 995.2
$ws.Range("J45").Value = 2604.4
$ws.Range("K45").Value = 995.2
$ws.Range("L45").Value = 2604.4
$ws.Range("M45").Value = -618.2
$ws.Range("N45").Value = -3358.4
$ws.Range("H74").Value = 2021260.2
$ws.Range("I74").Value = 1080.4681
$ws.Range("J74").Value = 13889816
$ws.Range("K74").Value = 1080.4681
$ws.Range("L74").Value = 13889816
$ws.Range("M74").Value = -206.4681
$ws.Range("N74").Value = -13891564
$ws.Range("H77").Value = 2021260.2
$ws.Range("I77").Value = 1080.4681
$ws.Range("J77").Value = 13889816
$ws.Range("K77").Value = 5402.3405
$ws.Range("L77").Value = 69449080
$ws.Range("M77").Value = -1034.3405
$ws.Range("N77").Value = -69457816
$ws.Range("H97").Value = 806.2857
$ws.Range("I97").Value = 489.81818
$ws.Range("J97").Value = 1966.6666
$ws.Range("K97").Value = 489.81818
$ws.Range("L97").Value = 1966.6666
$ws.Range("M97").Value = 6.181820000000016
$ws.Range("N97").Value = -2958.6666
$ws.Range("H102").Value = 1190
$ws.Range("I102").Value = 1190
$ws.Range("K102").Value = 1190
$ws.Range("M102").Value = 432
$ws.Range("H110").Value = 1825.6666
$ws.Range("I110").Value = 1318.1765
$ws.Range("J110").Value = 3982.5
$ws.Range("K110").Value = 1318.1765
$ws.Range("L110").Value = 3982.5
$ws.Range("M110").Value = 726.8235
$ws.Range("N110").Value = -8072.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1994.1305
$ws.Range("I134").Value = 2043.25
$ws.Range("J134").Value = 1666.6666
$ws.Range("K134").Value = 6129.75
$ws.Range("L134").Value = 4999.9998
$ws.Range("M134").Value = -3594.75
$ws.Range("N134").Value = -10069.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1653.7778
$ws.Range("I31").Value = 1362.1538
$ws.Range("J31").Value = 2127.6667
$ws.Range("K31").Value = 1362.1538
$ws.Range("L31").Value = 2127.6667
$ws.Range("M31").Value = -1067.1538
$ws.Range("N31").Value = -2717.6667
$ws.Range("H34").Value = 1653.7778
$ws.Range("I34").Value = 1362.1538
$ws.Range("J34").Value = 2127.6667
$ws.Range("K34").Value = 1362.1538
$ws.Range("L34").Value = 2127.6667
$ws.Range("M34").Value = -1160.1538
$ws.Range("N34").Value = -2531.6667
$ws.Range("H58").Value = 849.6949
$ws.Range("I58").Value = 751.67346
$ws.Range("J58").Value = 1330
$ws.Range("K58").Value = 751.67346
$ws.Range("L58").Value = 1330
$ws.Range("M58").Value = -548.67346
$ws.Range("N58").Value = -1736
$ws.Range("H132").Value = 1658.1818
$ws.Range("I132").Value = 1464.0667
$ws.Range("J132").Value = 3599.3333
$ws.Range("K132").Value = 4392.2001
$ws.Range("L132").Value = 10797.9999
$ws.Range("M132").Value = -1862.2001
$ws.Range("N132").Value = -15857.9999
$ws.Range("H134").Value = 1992.44
$ws.Range("I134").Value = 1008.3182
$ws.Range("J134").Value = 9209.333000000001
$ws.Range("K134").Value = 3024.9546
$ws.Range("L134").Value = 27627.999
$ws.Range("M134").Value = -489.9546
$ws.Range("N134").Value = -32697.999
$ws.Range("H136").Value = 849.6949
$ws.Range("I136").Value = 751.67346
$ws.Range("J136").Value = 1330
$ws.Range("K136").Value = 2255.02038
$ws.Range("L136").Value = 3990
$ws.Range("M136").Value = 294.9796200000001
$ws.Range("N136").Value = -9090

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 1847.091
$ws.Range("I108").Value = 1098.2858
$ws.Range("J108").Value = 3157.5
$ws.Range("K108").Value = 3294.8574
$ws.Range("L108").Value = 9472.5
$ws.Range("M108").Value = -414.8574000000003
$ws.Range("N108").Value = -15232.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 1033.8334
$ws.Range("I12").Value = 925.75
$ws.Range("J12").Value = 1250
$ws.Range("K12").Value = 925.75
$ws.Range("L12").Value = 1250
$ws.Range("M12").Value = -785.75
$ws.Range("N12").Value = -1530

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 7337.6665
$ws.Range("I11").Value = 2006
$ws.Range("J11").Value = 10003.5
$ws.Range("K11").Value = 2006
$ws.Range("L11").Value = 10003.5
$ws.Range("M11").Value = -1866
$ws.Range("N11").Value = -10283.5
$ws.Range("H132").Value = 1697.3
$ws.Range("I132").Value = 1278.0217
$ws.Range("J132").Value = 2135.6365
$ws.Range("K132").Value = 3834.0651
$ws.Range("L132").Value = 6406.9095
$ws.Range("M132").Value = -1304.0651
$ws.Range("N132").Value = -11466.9095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 1133
$ws.Range("I17").Value = 949.5
$ws.Range("J17").Value = 1500
$ws.Range("K17").Value = 949.5
$ws.Range("L17").Value = 1500
$ws.Range("M17").Value = -777.5
$ws.Range("N17").Value = -1844
$ws.Range("H19").Value = 4785.7144
$ws.Range("J19").Value = 4785.7144
$ws.Range("L19").Value = 4785.7144
$ws.Range("N19").Value = -5133.7144
$ws.Range("H136").Value = 4758.407
$ws.Range("I136").Value = 5648.65
$ws.Range("J136").Value = 2214.8572
$ws.Range("K136").Value = 16945.95
$ws.Range("L136").Value = 6644.571599999999
$ws.Range("M136").Value = -14395.95
$ws.Range("N136").Value = -11744.5716
